# Updates cryptos list values per latest scrape (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.779.63'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.266.46'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.15'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.95'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.484'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.43'
$ws.Range('E10').Value = '  +1.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.31'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.114'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.66'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.619.96'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.23'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.278.72'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.769'
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.682.19'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.53'
$ws.Range('E20').Value = '  +4.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0906'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.94'
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.09'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.42'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.92'
$ws.Range('E27').Value = '  +2.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.95'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.54'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.07'
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.32'
$ws.Range('E31').Value = '  +5.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.26'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0743'
$ws.Range('E35').Value = '  +1.25%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.93'
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.117'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.105'
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.90'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.022.48'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.30'
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.43'
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0279'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.33'
$ws.Range('E51').Value = '  +2.81%  '
